$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Target cluster changes from FAPs -> ECs, counts/values recomputed ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8330250000000001
$ws.Range("H2").Value = 2.499075
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 13.57958433333333
$ws.Range("N2").Value = 40.738753
$ws.Range("O2").Value = 0.2289698008477291
$ws.Range("P2").Value = 0.2289698008477291
$ws.Range("Q2").Value = 11.312133239275
$ws.Range("R2").Value = 101.809199153475
$ws.Range("S2").Value = 0.2289698008477291
$ws.Range("T2").Value = 0.2289698008477291

# --- Row 3: Target cluster changes from sCs -> FAPs, counts/values recomputed ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8330250000000001
$ws.Range("H3").Value = 2.499075
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.768727
$ws.Range("N3").Value = 59.306181
$ws.Range("O3").Value = 0.3333269541315948
$ws.Range("P3").Value = 0.3333269541315948
$ws.Range("Q3").Value = 16.467843809175
$ws.Range("R3").Value = 148.210594282575
$ws.Range("S3").Value = 0.3333269541315948
$ws.Range("T3").Value = 0.3333269541315948

# --- Row 4: brand-new row, Target cluster = sCs ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8330250000000001
$ws.Range("H4").Value = 2.499075
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 25.95900466666667
$ws.Range("N4").Value = 77.877014
$ws.Range("O4").Value = 0.4377032450206762
$ws.Range("P4").Value = 0.4377032450206762
$ws.Range("Q4").Value = 21.62449986245
$ws.Range("R4").Value = 194.62049876205
$ws.Range("S4").Value = 0.4377032450206762
$ws.Range("T4").Value = 0.4377032450206762
